# Update "想去人数" (number of people interested) counts for a handful of
# events in the "展览" sheet and the mirrored "全部类型" sheet.
#
#   F4  : 2172  -> 2173
#   F6  : 12716 -> 12719
#   F27 : 5176  -> 5179

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 2173
    $ws.Range("F6").Value = 12719
    $ws.Range("F27").Value = 5179
}
